# Commit automatique via PowerShell
# Sort the student list (rows 4-17) alphabetically by last name (column B),
# keeping the serial numbers in column A untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("B4:H17")
$key1 = $ws.Range("B4:B17")

$sortRange.Sort($key1, 1, $null, $null, 1, $null, 1, 0)
